# Scheduled runner update: refresh computed market-price / profit figures
# on a handful of rows across the per-job "Leve" profit sheets.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 19 on ALC
$ws_ALC.Range("H19").Value = 28425.066
$ws_ALC.Range("J19").Value = 52853.875
$ws_ALC.Range("L19").Value = 52853.875
$ws_ALC.Range("N19").Value = -53203.875

# Row 96 on ALC
$ws_ALC.Range("H96").Value = 1522.1428
$ws_ALC.Range("I96").Value = 1382.25
$ws_ALC.Range("J96").Value = 1708.6666
$ws_ALC.Range("K96").Value = 4146.75
$ws_ALC.Range("L96").Value = 5125.9998
$ws_ALC.Range("M96").Value = -2773.75
$ws_ALC.Range("N96").Value = -7871.9998

# Row 116 on ALC
$ws_ALC.Range("H116").Value = 2740.5652
$ws_ALC.Range("I116").Value = 2708.9402
$ws_ALC.Range("K116").Value = 2708.9402
$ws_ALC.Range("M116").Value = 733.0598

# Row 139 on ALC
$ws_ALC.Range("H139").Value = 98959.3
$ws_ALC.Range("J139").Value = 98959.3
$ws_ALC.Range("L139").Value = 98959.3
$ws_ALC.Range("N139").Value = -109239.3

$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 56 on ARM
$ws_ARM.Range("H56").Value = 24500
$ws_ARM.Range("I56").Value = 29000
$ws_ARM.Range("J56").Value = 20000
$ws_ARM.Range("K56").Value = 29000
$ws_ARM.Range("L56").Value = 20000
$ws_ARM.Range("M56").Value = -28258
$ws_ARM.Range("N56").Value = -21484

# Row 98 on ARM
$ws_ARM.Range("H98").Value = 118875
$ws_ARM.Range("J98").Value = 118875
$ws_ARM.Range("L98").Value = 118875
$ws_ARM.Range("N98").Value = -124865

# Row 102 on ARM
$ws_ARM.Range("H102").Value = 11218.964
$ws_ARM.Range("I102").Value = 4335.7393
$ws_ARM.Range("J102").Value = 42881.8
$ws_ARM.Range("K102").Value = 4335.7393
$ws_ARM.Range("L102").Value = 42881.8
$ws_ARM.Range("M102").Value = -2713.7393
$ws_ARM.Range("N102").Value = -46125.8

# Row 134 on ARM
$ws_ARM.Range("H134").Value = 96000
$ws_ARM.Range("J134").Value = 96000
$ws_ARM.Range("L134").Value = 96000
$ws_ARM.Range("N134").Value = -106140

# Row 135 on ARM
$ws_ARM.Range("H135").Value = 129499.5
$ws_ARM.Range("J135").Value = 129499.5
$ws_ARM.Range("L135").Value = 129499.5
$ws_ARM.Range("N135").Value = -139639.5

$ws_BSM = $wb.Worksheets.Item("BSM")
# Row 21 on BSM
$ws_BSM.Range("H21").Value = 0
$ws_BSM.Range("J21").Value = 0
$ws_BSM.Range("L21").Value = 0
$ws_BSM.Range("N21").ClearContents()

# Row 94 on BSM
$ws_BSM.Range("H94").Value = 5587.7437
$ws_BSM.Range("I94").Value = 4812.057
$ws_BSM.Range("K94").Value = 4812.057
$ws_BSM.Range("M94").Value = -4361.057

# Row 134 on BSM
$ws_BSM.Range("H134").Value = 11083.658
$ws_BSM.Range("I134").Value = 5889.154
$ws_BSM.Range("K134").Value = 17667.462
$ws_BSM.Range("M134").Value = -15132.462

$ws_CRP = $wb.Worksheets.Item("CRP")
# Row 74 on CRP
$ws_CRP.Range("H74").Value = 45000
$ws_CRP.Range("J74").Value = 45000
$ws_CRP.Range("L74").Value = 45000
$ws_CRP.Range("N74").Value = -46748

# Row 77 on CRP
$ws_CRP.Range("H77").Value = 45000
$ws_CRP.Range("J77").Value = 45000
$ws_CRP.Range("L77").Value = 135000
$ws_CRP.Range("N77").Value = -143736

# Row 132 on CRP
$ws_CRP.Range("H132").Value = 5023.108
$ws_CRP.Range("J132").Value = 12260.833
$ws_CRP.Range("L132").Value = 36782.499
$ws_CRP.Range("N132").Value = -41842.499

$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 63 on CUL
$ws_CUL.Range("H63").Value = 666
$ws_CUL.Range("I63").Value = 666
$ws_CUL.Range("K63").Value = 1998
$ws_CUL.Range("M63").Value = -1249

# Row 66 on CUL
$ws_CUL.Range("H66").Value = 666
$ws_CUL.Range("I66").Value = 666
$ws_CUL.Range("K66").Value = 5994
$ws_CUL.Range("M66").Value = -2250

# Row 128 on CUL
$ws_CUL.Range("H128").Value = 202778.6
$ws_CUL.Range("I128").Value = 202778.6
$ws_CUL.Range("K128").Value = 608335.8
$ws_CUL.Range("M128").Value = -603355.8

$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 14 on GSM
$ws_GSM.Range("H14").Value = 1125484.9
$ws_GSM.Range("I14").Value = 2250472.5
$ws_GSM.Range("J14").Value = 497.25
$ws_GSM.Range("K14").Value = 2250472.5
$ws_GSM.Range("L14").Value = 497.25
$ws_GSM.Range("M14").Value = -2250304.5
$ws_GSM.Range("N14").Value = -833.25

# Row 28 on GSM
$ws_GSM.Range("H28").Value = 24999.5
$ws_GSM.Range("J28").Value = 24999.5
$ws_GSM.Range("L28").Value = 24999.5
$ws_GSM.Range("N28").Value = -25383.5

# Row 39 on GSM
$ws_GSM.Range("H39").Value = 15253.571
$ws_GSM.Range("J39").Value = 15253.571
$ws_GSM.Range("L39").Value = 15253.571
$ws_GSM.Range("N39").Value = -16317.571

# Row 69 on GSM
$ws_GSM.Range("H69").Value = 46275.715
$ws_GSM.Range("J69").Value = 46275.715
$ws_GSM.Range("L69").Value = 46275.715
$ws_GSM.Range("N69").Value = -47773.715

# Row 72 on GSM
$ws_GSM.Range("H72").Value = 46275.715
$ws_GSM.Range("J72").Value = 46275.715
$ws_GSM.Range("L72").Value = 138827.145
$ws_GSM.Range("N72").Value = -146315.145

# Row 75 on GSM
$ws_GSM.Range("H75").Value = 100000
$ws_GSM.Range("J75").Value = 100000
$ws_GSM.Range("L75").Value = 100000
$ws_GSM.Range("N75").Value = -101748

# Row 78 on GSM
$ws_GSM.Range("H78").Value = 100000
$ws_GSM.Range("J78").Value = 100000
$ws_GSM.Range("L78").Value = 300000
$ws_GSM.Range("N78").Value = -308736

# Row 80 on GSM
$ws_GSM.Range("H80").Value = 15926.55
$ws_GSM.Range("I80").Value = 8478.416999999999
$ws_GSM.Range("J80").Value = 27098.75
$ws_GSM.Range("K80").Value = 8478.416999999999
$ws_GSM.Range("L80").Value = 27098.75
$ws_GSM.Range("M80").Value = -7480.416999999999
$ws_GSM.Range("N80").Value = -29094.75

# Row 83 on GSM
$ws_GSM.Range("H83").Value = 15926.55
$ws_GSM.Range("I83").Value = 8478.416999999999
$ws_GSM.Range("J83").Value = 27098.75
$ws_GSM.Range("K83").Value = 42392.085
$ws_GSM.Range("L83").Value = 135493.75
$ws_GSM.Range("M83").Value = -37400.085
$ws_GSM.Range("N83").Value = -145477.75

# Row 93 on GSM
$ws_GSM.Range("H93").Value = 46689.86
$ws_GSM.Range("J93").Value = 46689.86
$ws_GSM.Range("L93").Value = 46689.86
$ws_GSM.Range("N93").Value = -50433.86

# Row 95 on GSM
$ws_GSM.Range("H95").Value = 30000
$ws_GSM.Range("J95").Value = 30000
$ws_GSM.Range("L95").Value = 30000
$ws_GSM.Range("N95").Value = -35492

# Row 135 on GSM
$ws_GSM.Range("H135").Value = 154098.88
$ws_GSM.Range("J135").Value = 154098.88
$ws_GSM.Range("L135").Value = 154098.88
$ws_GSM.Range("N135").Value = -164238.88

$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 22 on LTW
$ws_LTW.Range("H22").Value = 7380
$ws_LTW.Range("I22").Value = 1200
$ws_LTW.Range("J22").Value = 8925
$ws_LTW.Range("K22").Value = 1200
$ws_LTW.Range("L22").Value = 8925
$ws_LTW.Range("M22").Value = -905
$ws_LTW.Range("N22").Value = -9515

# Row 27 on LTW
$ws_LTW.Range("H27").Value = 7380
$ws_LTW.Range("I27").Value = 1200
$ws_LTW.Range("J27").Value = 8925
$ws_LTW.Range("K27").Value = 1200
$ws_LTW.Range("L27").Value = 8925
$ws_LTW.Range("M27").Value = -1093
$ws_LTW.Range("N27").Value = -9139

# Row 68 on LTW
$ws_LTW.Range("H68").Value = 18286.334
$ws_LTW.Range("I68").Value = 0
$ws_LTW.Range("J68").Value = 18286.334
$ws_LTW.Range("K68").Value = 0
$ws_LTW.Range("L68").Value = 18286.334
$ws_LTW.Range("M68").ClearContents()
$ws_LTW.Range("N68").Value = -19784.334

# Row 71 on LTW
$ws_LTW.Range("H71").Value = 18286.334
$ws_LTW.Range("I71").Value = 0
$ws_LTW.Range("J71").Value = 18286.334
$ws_LTW.Range("K71").Value = 0
$ws_LTW.Range("L71").Value = 91431.67
$ws_LTW.Range("M71").ClearContents()
$ws_LTW.Range("N71").Value = -98919.67

# Row 87 on LTW
$ws_LTW.Range("H87").Value = 100189
$ws_LTW.Range("J87").Value = 100189
$ws_LTW.Range("L87").Value = 100189
$ws_LTW.Range("N87").Value = -102435

# Row 90 on LTW
$ws_LTW.Range("H90").Value = 100189
$ws_LTW.Range("J90").Value = 100189
$ws_LTW.Range("L90").Value = 300567
$ws_LTW.Range("N90").Value = -311799

$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 136 on WVR
$ws_WVR.Range("H136").Value = 16235.529
$ws_WVR.Range("I136").Value = 4001.3333
$ws_WVR.Range("J136").Value = 18857.143
$ws_WVR.Range("K136").Value = 12003.9999
$ws_WVR.Range("L136").Value = 56571.429
$ws_WVR.Range("M136").Value = -9453.999899999999
$ws_WVR.Range("N136").Value = -61671.429
